$d = $word.ActiveDocument
$x = [char]0xD7

$replacements = @(
    @{old = "2025-02-12 Wednesday"; new = "2025-02-13 Thursday"},
    @{old = "63$($x)34="; new = "85$($x)23="},
    @{old = "32$($x)62="; new = "39$($x)98="},
    @{old = "34$($x)35="; new = "33$($x)84="},
    @{old = "81$($x)35="; new = "94$($x)92="},
    @{old = "56$($x)55="; new = "75$($x)65="},
    @{old = "30$($x)57="; new = "92$($x)65="},
    @{old = "15$($x)41="; new = "12$($x)14="},
    @{old = "73$($x)49="; new = "54$($x)34="},
    @{old = "46$($x)95="; new = "56$($x)52="},
    @{old = "14$($x)33="; new = "20$($x)14="},
    @{old = "24$($x)73="; new = "88$($x)14="},
    @{old = "69$($x)40="; new = "41$($x)70="},
    @{old = "86$($x)44="; new = "89$($x)42="},
    @{old = "94$($x)19="; new = "22$($x)55="},
    @{old = "26$($x)84="; new = "78$($x)72="},
    @{old = "69$($x)53="; new = "84$($x)62="},
    @{old = "96$($x)27="; new = "13$($x)89="},
    @{old = "83$($x)36="; new = "84$($x)57="},
    @{old = "12$($x)54="; new = "20$($x)35="},
    @{old = "36$($x)11="; new = "43$($x)33="},
    @{old = "96$($x)40="; new = "91$($x)65="},
    @{old = "79$($x)73="; new = "16$($x)40="},
    @{old = "45$($x)79="; new = "56$($x)26="},
    @{old = "39$($x)83="; new = "28$($x)49="},
    @{old = "29$($x)76="; new = "89$($x)57="}
)

foreach ($rep in $replacements) {
    $d.Content.Find.Execute($rep.old, $true, $false, $false, $false, $false, $true, 1, $false, $rep.new, 2) | Out-Null
}
